$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN) to hold "Week_Start_Date"
$ws.Columns("B").Insert()

# New header for the inserted column
$ws.Range("B1").Value = "Week_Start_Date"

# Week start dates for each of the 16 forecast weeks (rows 2-17)
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    # Leading apostrophe forces Excel to store the value as literal text
    # instead of auto-converting the date-like string to a date serial.
    $ws.Range("B$row").Value = "'" + $weekStartDates[$i]
}

# Correct the Week labels in column A: "W01".."W09" -> "W1".."W9"
# (W10..W16 are already correct and unchanged)
for ($i = 1; $i -le 9; $i++) {
    $row = $i + 1
    $ws.Range("A$row").Value = "W$i"
}
